# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'43.085.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "`'2.308.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "`'300.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "`'98.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").Value = "`'36.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "`'0.0793"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "`'18.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").Value = "`'2.667.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "`'2.303.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "`'0.782"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "`'43.017.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "`'12.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.35%  "
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "`'6.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").Value = "`'68.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "`'240.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "`'2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "`'2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "`'25.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("D29").Value = "`'165.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").Value = "`'33.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "`'4.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").Value = "`'17.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.47%  "
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D38").Value = "`'0.0688"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  -1.76%  "
$ws.Range("D43").Value = "`'2.021.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").Value = "`'0.0283"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").Value = "`'17.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").Value = "`'2.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "`'54.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "`'2.536.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
